$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.669.75'
$ws.Range("E2").Value = '  -2.79%  '

$ws.Range("D3").Value = '1.982.66'
$ws.Range("E3").Value = '  -3.73%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.634'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.57%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.11'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.86%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '58.98'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.05%  '

$ws.Range("E10").Value = '  -0.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0738'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.92%  '

$ws.Range("E12").Value = '  -2.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.956'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.78%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.72'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.56%  '

$ws.Range("D15").Value = '2.271.75'

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.115.88'
$ws.Range("E16").Value = '  +2.47%  '

$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.90%  '

$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.91'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +10.36%  '

$ws.Range("D19").Value = '35.563.16'
$ws.Range("E19").Value = '  -2.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.46'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.98%  '

$ws.Range("D21").Value = '0.0₃0847'
$ws.Range("E21").Value = '  -2.34%  '

$ws.Range("E22").Value = '  -0.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.21%  '

$ws.Range("E24").Value = '  +0.11%  '

$ws.Range("E25").Value = '  +19.95%  '

$ws.Range("E26").Value = '  -4.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.62%  '

$ws.Range("E30").Value = '  -2.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.90'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.24%  '

$ws.Range("E32").Value = '  -7.15%  '

$ws.Range("E33").Value = '  +12.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0597'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.39%  '

$ws.Range("E35").Value = '  +9.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.38'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.30%  '

$ws.Range("E37").Value = '  -0.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.99%  '

$ws.Range("E39").Value = '  +9.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.23'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.53%  '

$ws.Range("E41").Value = '  -1.19%  '

$ws.Range("E42").Value = '  -0.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '93.83'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.07%  '

$ws.Range("E44").Value = '  -1.39%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0911'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.39%  '

$ws.Range("D48").Value = '1.372.86'
$ws.Range("E48").Value = '  -2.89%  '

$ws.Range("E49").Value = '  -0.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '46.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.27%  '
